$wb = $excel.ActiveWorkbook

# --- Rename the "Consumer_Sectors" sheet to "Consumer_Firm_Sectors" ---
$wsConsumer = $wb.Worksheets.Item("Consumer_Sectors")
$wsConsumer.Name = "Consumer_Firm_Sectors"

# --- Update values on Main_Loop_Parameters ---
$wsMain = $wb.Worksheets.Item("Main_Loop_Parameters")
$wsMain.Range("B3").Value = 60
$wsMain.Range("B4").Value = 10000
$wsMain.Range("B5").Value = 400
$wsMain.Range("B6").Value = 50
$wsMain.Range("B7").Formula = "=B5+B6"

# --- Update selection on Consumer_Firm_Sectors (was the active tab, now just has a new selection) ---
$wsConsumer.Activate()
$wsConsumer.Range("C8").Select()

# --- Update selection on Initialization_Parameters (not the active tab, but keeps a stored selection) ---
$wsInit = $wb.Worksheets.Item("Initialization_Parameters")
$wsInit.Activate()
$wsInit.Range("B4").Select()

# --- Main_Loop_Parameters becomes the active tab with B4 selected ---
$wsMain.Activate()
$wsMain.Range("B4").Select()

$wb.Save()
